$wb = $excel.ActiveWorkbook
$ds = $wb.Worksheets.Item("data")

# --- Add the new "metadata" sheet, positioned after "data" ---
$md = $wb.Worksheets.Add($null, $ds)
$md.Name = "metadata"

# Header row (bold, bordered, centered — matching the "data" sheet's header style)
$md.Range("B1").Value = "data_name"
$md.Range("C1").Value = "data_id"
$md.Range("D1").Value = "data_version"
$md.Range("E1").Value = "data_version_created"
$md.Range("F1").Value = "panel_query_time"
$md.Range("G1").Value = "panel_get_request"

$md.Range("B1:G1").Font.Bold = $true
$md.Range("B1:G1").Borders.LineStyle = 1
$md.Range("B1:G1").HorizontalAlignment = -4108
$md.Range("B1:G1").VerticalAlignment = -4160

# Index cell A2 (bold, bordered, centered — matching column A style on "data")
$md.Range("A2").Value = 0
$md.Range("A2").Font.Bold = $true
$md.Range("A2").Borders.LineStyle = 1
$md.Range("A2").HorizontalAlignment = -4108
$md.Range("A2").VerticalAlignment = -4160

# Data row
$md.Range("B2").Value = "White matter disorders - adult onset"
$md.Range("C2").Value = 579
$md.Range("D2").Value = "'1.30"
$md.Range("E2").Value = "2021-08-31T15:05:14.341243Z"
$md.Range("F2").Value = "2021-10-05 14:23:08.349782"
$md.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/579/?format=json"

# Update time_taken (column F) timestamps on the "data" sheet
$ds.Range("F2").Value = "2021-10-05 14:23:08.353082"
$ds.Range("F3").Value = "2021-10-05 14:23:08.353090"
$ds.Range("F4").Value = "2021-10-05 14:23:08.353093"
$ds.Range("F5").Value = "2021-10-05 14:23:08.353096"
$ds.Range("F6").Value = "2021-10-05 14:23:08.353099"
$ds.Range("F7").Value = "2021-10-05 14:23:08.353102"
$ds.Range("F8").Value = "2021-10-05 14:23:08.353105"
$ds.Range("F9").Value = "2021-10-05 14:23:08.353107"
$ds.Range("F10").Value = "2021-10-05 14:23:08.353110"
$ds.Range("F11").Value = "2021-10-05 14:23:08.353113"
$ds.Range("F12").Value = "2021-10-05 14:23:08.353116"
$ds.Range("F13").Value = "2021-10-05 14:23:08.353118"
$ds.Range("F14").Value = "2021-10-05 14:23:08.353121"
$ds.Range("F15").Value = "2021-10-05 14:23:08.353124"
$ds.Range("F16").Value = "2021-10-05 14:23:08.353127"
$ds.Range("F17").Value = "2021-10-05 14:23:08.353129"
$ds.Range("F18").Value = "2021-10-05 14:23:08.353132"
$ds.Range("F19").Value = "2021-10-05 14:23:08.353135"
$ds.Range("F20").Value = "2021-10-05 14:23:08.353137"
$ds.Range("F21").Value = "2021-10-05 14:23:08.353140"
$ds.Range("F22").Value = "2021-10-05 14:23:08.353143"
$ds.Range("F23").Value = "2021-10-05 14:23:08.353145"
$ds.Range("F24").Value = "2021-10-05 14:23:08.353148"
$ds.Range("F25").Value = "2021-10-05 14:23:08.353151"
$ds.Range("F26").Value = "2021-10-05 14:23:08.353154"
$ds.Range("F27").Value = "2021-10-05 14:23:08.353157"
$ds.Range("F28").Value = "2021-10-05 14:23:08.353159"
$ds.Range("F29").Value = "2021-10-05 14:23:08.353162"
$ds.Range("F30").Value = "2021-10-05 14:23:08.353165"
$ds.Range("F31").Value = "2021-10-05 14:23:08.353167"
$ds.Range("F32").Value = "2021-10-05 14:23:08.353170"
$ds.Range("F33").Value = "2021-10-05 14:23:08.353172"
$ds.Range("F34").Value = "2021-10-05 14:23:08.353175"
$ds.Range("F35").Value = "2021-10-05 14:23:08.353178"
$ds.Range("F36").Value = "2021-10-05 14:23:08.353181"
$ds.Range("F37").Value = "2021-10-05 14:23:08.353184"
$ds.Range("F38").Value = "2021-10-05 14:23:08.353186"
$ds.Range("F39").Value = "2021-10-05 14:23:08.353189"
$ds.Range("F40").Value = "2021-10-05 14:23:08.353192"
$ds.Range("F41").Value = "2021-10-05 14:23:08.353194"
$ds.Range("F42").Value = "2021-10-05 14:23:08.353197"
$ds.Range("F43").Value = "2021-10-05 14:23:08.353200"
$ds.Range("F44").Value = "2021-10-05 14:23:08.353203"
$ds.Range("F45").Value = "2021-10-05 14:23:08.353205"
$ds.Range("F46").Value = "2021-10-05 14:23:08.353208"
$ds.Range("F47").Value = "2021-10-05 14:23:08.353210"
$ds.Range("F48").Value = "2021-10-05 14:23:08.353213"
$ds.Range("F49").Value = "2021-10-05 14:23:08.353215"
$ds.Range("F50").Value = "2021-10-05 14:23:08.353218"
$ds.Range("F51").Value = "2021-10-05 14:23:08.353220"
$ds.Range("F52").Value = "2021-10-05 14:23:08.353223"
$ds.Range("F53").Value = "2021-10-05 14:23:08.353226"
$ds.Range("F54").Value = "2021-10-05 14:23:08.353229"
$ds.Range("F55").Value = "2021-10-05 14:23:08.353231"
$ds.Range("F56").Value = "2021-10-05 14:23:08.353234"
$ds.Range("F57").Value = "2021-10-05 14:23:08.353237"
$ds.Range("F58").Value = "2021-10-05 14:23:08.353239"
$ds.Range("F59").Value = "2021-10-05 14:23:08.353242"
$ds.Range("F60").Value = "2021-10-05 14:23:08.353244"
$ds.Range("F61").Value = "2021-10-05 14:23:08.353247"
$ds.Range("F62").Value = "2021-10-05 14:23:08.353250"
$ds.Range("F63").Value = "2021-10-05 14:23:08.353253"
$ds.Range("F64").Value = "2021-10-05 14:23:08.353255"
$ds.Range("F65").Value = "2021-10-05 14:23:08.353258"
$ds.Range("F66").Value = "2021-10-05 14:23:08.353262"
$ds.Range("F67").Value = "2021-10-05 14:23:08.353265"
$ds.Range("F68").Value = "2021-10-05 14:23:08.353268"
$ds.Range("F69").Value = "2021-10-05 14:23:08.353270"
$ds.Range("F70").Value = "2021-10-05 14:23:08.353273"
$ds.Range("F71").Value = "2021-10-05 14:23:08.353275"
$ds.Range("F72").Value = "2021-10-05 14:23:08.353278"
$ds.Range("F73").Value = "2021-10-05 14:23:08.353281"
$ds.Range("F74").Value = "2021-10-05 14:23:08.353284"
$ds.Range("F75").Value = "2021-10-05 14:23:08.353286"
$ds.Range("F76").Value = "2021-10-05 14:23:08.353289"
$ds.Range("F77").Value = "2021-10-05 14:23:08.353292"
$ds.Range("F78").Value = "2021-10-05 14:23:08.353296"
$ds.Range("F79").Value = "2021-10-05 14:23:08.353300"
$ds.Range("F80").Value = "2021-10-05 14:23:08.353303"
$ds.Range("F81").Value = "2021-10-05 14:23:08.353305"
$ds.Range("F82").Value = "2021-10-05 14:23:08.353308"
$ds.Range("F83").Value = "2021-10-05 14:23:08.353311"
$ds.Range("F84").Value = "2021-10-05 14:23:08.353313"
$ds.Range("F85").Value = "2021-10-05 14:23:08.353316"
$ds.Range("F86").Value = "2021-10-05 14:23:08.353318"
$ds.Range("F87").Value = "2021-10-05 14:23:08.353321"
$ds.Range("F88").Value = "2021-10-05 14:23:08.353324"
$ds.Range("F89").Value = "2021-10-05 14:23:08.353326"

# Restore "data" as the active sheet (unchanged in the source bookViews)
$ds.Activate()
